# CLNEBag.xlsx edit
# - Adds a new "totalScore" column right after "Date"
# - Drops the sample data row (row 2), leaving only the header row
# - Removes the now-unused "Bag" header / value and the date formatting
#   that used to live on column A

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sample data row - only the header row remains.
$ws.Rows.Item(2).Delete()

# Insert a new column for "totalScore" right after the Date column.
$ws.Columns.Item(2).Insert()
$ws.Cells.Item(1, 2).Value = "totalScore"

# The Date column no longer carries a date cell, so drop its number
# format (and the new totalScore column shouldn't inherit it either).
$ws.Columns.Item(1).ClearFormats()
$ws.Columns.Item(2).ClearFormats()

# Resize the Date / totalScore columns to fit their header text; the
# remaining columns keep the widths they already had.
$ws.Columns.Item(1).ColumnWidth = 4.3
$ws.Columns.Item(2).ColumnWidth = 9.15
